# cardinality.pptx edit script
# - Delete the old slide 15 ("Rationals are countable") - superseded by the
#   "Reals are ..." slide that follows it.
# - Fix the title of the slide that takes its place (formerly slide 16,
#   "Reals are uncountable") to read "Reals are countable".
# - Mark slides 10-14 as hidden in the slide show (show="0").

$p = $ppt.ActivePresentation

# Remove the old "Rationals are countable" slide (slide 15).
$p.Slides.Item(15).Delete()

# The slide that used to be #16 ("Reals are uncountable") is now #15.
$s15 = $p.Slides.Item(15)
$title = $s15.Shapes.Item(2)
$tr = $title.TextFrame.TextRange
$tr.Characters(11, 2).Delete()

# Hide slides 10-14 during the slide show.
foreach ($i in 10,11,12,13,14) {
    $p.Slides.Item($i).SlideShowTransition.Hidden = $true
}
